# Update cryptocurrency price/volume data per the source refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.982.72'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '1.744.83'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('D4').Value = "'1.000"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'249.68"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.58%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = "'0.5145"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.03%  '
$ws.Range('D8').Value = "'0.2760"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.52%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = '1.742.14'
$ws.Range('E10').Value = '  -0.93%  '
$ws.Range('D11').Value = "'0.07239"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('D12').Value = "'15.22"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('D13').Value = "'0.6497"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('D14').Value = "'4.637"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = "'77.77"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.62%  '
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').Value = "'0.9999"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = '26.012.69'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('E19').Value = '  +1.69%  '
$ws.Range('D20').Value = "'0.000006818"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.59%  '
$ws.Range('D21').Value = '1.964.34'
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('D22').Value = "'4.301"
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = "'8.681"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.69%  '
$ws.Range('D24').Value = "'5.369"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.13%  '
$ws.Range('D25').Value = "'135.71"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.14%  '
$ws.Range('D26').Value = "'1.510"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('D28').Value = "'1.787"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.63%  '
$ws.Range('D29').Value = "'106.18"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.79%  '
$ws.Range('D30').Value = "'3.949"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.28%  '
$ws.Range('D31').Value = "'0.08251"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.23%  '
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('D33').Value = "'0.04677"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.63%  '
$ws.Range('D34').Value = "'2.656"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('D35').Value = "'1.001"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').Value = "'0.6255"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.12%  '
$ws.Range('D37').Value = "'2.726"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.32%  '
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('D39').Value = "'1.931"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.92%  '
$ws.Range('D40').Value = "'0.9998"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').Value = "'100.57"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('D42').Value = "'0.3885"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.55%  '
$ws.Range('D43').Value = "'0.7565"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.70%  '
$ws.Range('D44').Value = "'5.026"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.72%  '
$ws.Range('D45').Value = "'6.358"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('E46').Value = '  -0.60%  '
$ws.Range('E47').Value = '  +2.85%  '
$ws.Range('D48').Value = "'0.05237"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('D49').Value = "'30.71"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('D50').Value = "'7.595"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('D51').Value = "'0.3439"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.89%  '
